$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 750 (the "17世紀に発明されて以来…" post) and shift subsequent rows up
$ws.Rows.Item(750).Delete()
